$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.112349152565002
$ws.Range("B1").Value = 1.874949336051941
$ws.Range("C1").Value = 5.44464635848999
$ws.Range("D1").Value = 0.6820975542068481
$ws.Range("E1").Value = 0.7080945372581482
